$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used by a few of the resource names (matches source data).
$nbsp = [char]0x00A0

# Final (post-edit) state of the resource table, sorted by resource id with
# each "Rn" row immediately followed by its "RnX" waste row.
$data = @(
    @("R1",   1,    "Population"),
    @("R2",   2,    "MetallicElements"),
    @("R2X",  -0.1, "MetallicElementsWaste"),
    @("R3",   2,    "Timber"),
    @("R3X",  -0.1, "TimberWaste"),
    @("R4",   2,    "Land"),
    @("R4X",  -0.1, "LandWaste"),
    @("R5",   2,    "RenewableEnergyCapacity"),
    @("R5X",  -0.1, "RECWaste"),
    @("R6",   1,    "Fossil energy capacity"),
    @("R6X",  -0.1, "FECWaste"),
    @("R7",   3,    "Water"),
    @("R7X",  -0.1, "Water Waste"),
    @("R8",   5,    "Farm$nbsp"),
    @("R8X",  -0.2, "FarmWaste$nbsp"),
    @("R20",  5,    "Military"),
    @("R20X", -0.2, "WasteMilitary"),
    @("R21",  5,    "MetallicAlloys"),
    @("R21X", -0.1, "WasteMetallicAlloys"),
    @("R22",  5,    "Housing"),
    @("R22X", -0.1, "WasteHousing"),
    @("R23",  5,    "Food"),
    @("R23X", -0.1, "WasteFood"),
    @("R24",  8,    "PreparedFossilEnergy$nbsp"),
    @("R24X", -0.3, "WastePreparedFossilEnergy$nbsp"),
    @("R25",  10,   "Electronics"),
    @("R25X", -0.2, "WasteElectronics"),
    @("R26",  10,   "PreparedRenewableEnergy"),
    @("R26X", -0.2, "WastePreparedRenewableEnergy$nbsp")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# Column C now holds longer waste-resource names; widen it to fit ("bestFit").
$ws.Columns.Item(3).EntireColumn.AutoFit()

# Record where the cursor ended up after the edit.
$ws.Range("C25").Select()
